# Update countries & provincias Spain
# Applies the 28 Abril 2020 16:52 data refresh to the "Pais" sheet:
#  - bump the "Datos actualizados..." timestamp
#  - refresh several countries' case counters
#  - Bulgaria now sorts above Cuba (Bulgaria's refreshed row takes row 81,
#    Cuba's untouched row slides down to row 82)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 16:52"

# --- Estados Unidos (row 4) --------------------------------------------
$ws.Range("B4").Value = 1012855
$ws.Range("C4").Value = 2499
$ws.Range("E4").Value = 816420
$ws.Range("G4").Value = 219
$ws.Range("H4").Value = 57016

# --- Alemania (row 8) ---------------------------------------------------
$ws.Range("B8").Value = 159038
$ws.Range("C8").Value = 280
$ws.Range("E8").Value = 35477
$ws.Range("G8").Value = 35
$ws.Range("H8").Value = 6161

# --- Paises Bajos (row 17) ----------------------------------------------
$ws.Range("F17").Value = 861

# --- Bulgaria / Cuba swap (rows 81-82) ----------------------------------
# Row 81 becomes Bulgaria with refreshed numbers; row 82 becomes Cuba with
# its previous (unchanged) numbers.
$ws.Range("A81").Value = "Bulgaria"
$ws.Range("B81").Value = 1399
$ws.Range("C81").Value = 36
$ws.Range("D81").Value = 222
$ws.Range("E81").Value = 1119
$ws.Range("F81").Value = 39
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 58

$ws.Range("A82").Value = "Cuba"
$ws.Range("B82").Value = 1389
$ws.Range("C82").Value = 0
$ws.Range("D82").Value = 525
$ws.Range("E82").Value = 808
$ws.Range("F82").Value = 12
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 56

# --- Sri Lanka (row 105) -------------------------------------------------
$ws.Range("B105").Value = 599
$ws.Range("C105").Value = 11
$ws.Range("E105").Value = 458

# --- Jordania (row 113) --------------------------------------------------
$ws.Range("D113").Value = 348
$ws.Range("E113").Value = 93
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 8

# --- Mauricio (row 121) ---------------------------------------------------
$ws.Range("D121").Value = 303
$ws.Range("E121").Value = 21
